$d = $word.ActiveDocument
$sec = $d.Sections.First
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        Write-Output "Header $i InlineShapes: $($hdr.Range.InlineShapes.Count)"
        for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
            $shp = $hdr.Range.InlineShapes.Item($j)
            Write-Output "  Shape $j Title=[$($shp.Title)] Alt=[$($shp.AlternativeText)]"
        }
    }
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        Write-Output "Footer $i InlineShapes: $($ftr.Range.InlineShapes.Count)"
        for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
            $shp = $ftr.Range.InlineShapes.Item($j)
            Write-Output "  Shape $j Title=[$($shp.Title)] Alt=[$($shp.AlternativeText)]"
        }
    }
}
